$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update group index column + value column for each of the 5 areas ---
# Area1: A/B, Area2: E/F, Area3: I/J, Area4: M/N, Area5: Q/R

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 49272
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 21355
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 28671
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 19394
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 49427

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2634
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 4491
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3120
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1327
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1400

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1219
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 7724
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 12159
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 3531
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 3936

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2288
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 13469
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 10505
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 35227
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 5446

# --- Move the border-styled blank cells from column I to column H on rows 12-15 ---
# Use A11 (which already carries the target "bordered / bold / centered" style)
# as a format donor so the same cellXf (style index 1) gets reused instead of a
# brand new style being created.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I12").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("I15").ClearContents()

# --- Selection change ---
$ws.Range("L22").Select()

$wb.Save()
